$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.360.65'
$ws.Range("D3").Value = '1.870.93'
$ws.Range("E3").Value = '  -0.38%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '''238.00'
$ws.Range("E5").Value = '  +0.68%  '
$ws.Range("D6").Value = '''1.001'
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("D7").Value = '''0.4824'
$ws.Range("E7").Value = '  -0.57%  '
$ws.Range("D8").Value = '1.870.58'
$ws.Range("E8").Value = '  -0.40%  '
$ws.Range("D9").Value = '''0.2816'
$ws.Range("E9").Value = '  -2.15%  '
$ws.Range("D10").Value = '''0.06511'
$ws.Range("E10").Value = '  -1.09%  '
$ws.Range("D11").Value = '1.875.58'
$ws.Range("E11").Value = '  -0.17%  '
$ws.Range("D12").Value = '''0.07487'
$ws.Range("E12").Value = '  +2.43%  '
$ws.Range("D13").Value = '''16.47'
$ws.Range("E13").Value = '  -1.69%  '
$ws.Range("D14").Value = '''5.074'
$ws.Range("E14").Value = '  -1.49%  '
$ws.Range("D15").Value = '''87.92'
$ws.Range("E15").Value = '  +0.85%  '
$ws.Range("D16").Value = '''0.6562'
$ws.Range("E16").Value = '  +0.18%  '
$ws.Range("D17").Value = '30.326.84'
$ws.Range("E17").Value = '  +0.34%  '
$ws.Range("D18").Value = '''13.25'
$ws.Range("E18").Value = '  -0.70%  '
$ws.Range("D19").Value = '''1.001'
$ws.Range("E19").Value = '  +0.07%  '
$ws.Range("D20").Value = '''0.000007604'
$ws.Range("E20").Value = '  -1.62%  '
$ws.Range("D21").Value = '2.113.31'
$ws.Range("E21").Value = '  -0.71%  '
$ws.Range("D22").Value = '''222.62'
$ws.Range("E22").Value = '  +14.46%  '
$ws.Range("D23").Value = '''5.293'
$ws.Range("E23").Value = '  -0.46%  '
$ws.Range("D24").Value = '''1.001'
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("D25").Value = '''6.171'
$ws.Range("E25").Value = '  +0.82%  '
$ws.Range("D26").Value = '''9.262'
$ws.Range("E26").Value = '  -0.27%  '
$ws.Range("D27").Value = '''166.51'
$ws.Range("E27").Value = '  +3.42%  '
$ws.Range("D28").Value = '''18.69'
$ws.Range("E28").Value = '  +3.87%  '
$ws.Range("D29").Value = '''1.980'
$ws.Range("E29").Value = '  +3.41%  '
$ws.Range("D30").Value = '''1.461'
$ws.Range("E30").Value = '  +1.57%  '
$ws.Range("D31").Value = '''0.09371'
$ws.Range("E31").Value = '  +2.73%  '
$ws.Range("D32").Value = '''4.317'
$ws.Range("E32").Value = '  +1.05%  '
$ws.Range("D33").Value = '''4.010'
$ws.Range("E33").Value = '  -0.56%  '
$ws.Range("D34").Value = '''0.05041'
$ws.Range("E34").Value = '  -0.92%  '
$ws.Range("D35").Value = '''1.210'
$ws.Range("E35").Value = '  +10.35%  '
$ws.Range("D36").Value = '''0.7453'
$ws.Range("E36").Value = '  +3.73%  '
$ws.Range("D37").Value = '''2.711'
$ws.Range("E37").Value = '  +0.48%  '
$ws.Range("D38").Value = '''0.01830'
$ws.Range("E38").Value = '  +1.89%  '
$ws.Range("D39").Value = '''2.625'
$ws.Range("E39").Value = '  -0.62%  '
$ws.Range("D40").Value = '''2.084'
$ws.Range("E40").Value = '  +2.01%  '
$ws.Range("D41").Value = '''0.9082'
$ws.Range("E41").Value = '  -1.10%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '''5.920'
$ws.Range("E42").Value = '  +2.06%  '
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").Value = '''106.49'
$ws.Range("E43").Value = '  +0.27%  '
$ws.Range("D44").Value = '''0.4279'
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("E45").Value = '  +0.28%  '
$ws.Range("D46").Value = '''7.422'
$ws.Range("E46").Value = '  +0.11%  '
$ws.Range("D47").Value = '''0.1298'
$ws.Range("E47").Value = '  -1.26%  '
$ws.Range("D48").Value = '''64.29'
$ws.Range("E48").Value = '  -1.60%  '
$ws.Range("D49").Value = '''9.047'
$ws.Range("E49").Value = '  +1.05%  '
$ws.Range("D50").Value = '''1.487'
$ws.Range("E50").Value = '  +9.11%  '
$ws.Range("D51").Value = '''33.99'
$ws.Range("E51").Value = '  +0.37%  '
